$d = $word.ActiveDocument

# 1. Merge "Name: Sandesh Varma ... Date: " runs (text content unchanged, just recombined)
$findText1 = "Name: Sandesh Varma                                                                      Date: "
$d.Content.Find.Execute($findText1, $true, $false, $false, $false, $false, $true, 1, $false, $findText1, 2) | Out-Null

# 2. Merge "08/09" + "/2025" into a single run "08/09/2025"
$d.Content.Find.Execute("08/09/2025", $true, $false, $false, $false, $false, $true, 1, $false, "08/09/2025", 2) | Out-Null

# 3. Merge "Experiment No: " + "5" into a single run "Experiment No: 5"
$d.Content.Find.Execute("Experiment No: 5", $true, $false, $false, $false, $false, $true, 1, $false, "Experiment No: 5", 2) | Out-Null

# 4. Merge the two underscore runs into a single longer underscore run
$oldUnderscore = "__________________________________________________" + "_________________"
$newUnderscore = "___________________________________________________________________"
$d.Content.Find.Execute($oldUnderscore, $true, $false, $false, $false, $false, $true, 1, $false, $newUnderscore, 2) | Out-Null

# 5. Move the "_GoBack" bookmark from its paragraph near the end to the very start of the
#    document (immediately before the first run of the first paragraph).
#    Inserting at absolute position 0 has quirky bookmark-placement behaviour in this
#    runtime, so we insert a temporary placeholder character, anchor the bookmark right
#    after it, then remove the placeholder again.
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")
$anchorRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $anchorRange) | Out-Null
$d.Range(0, 1).Delete() | Out-Null

# 6. Add page borders to the section properties
$section = $d.Sections.Item(1)
$borders = $section.Borders
$borders.Item(-4).LineStyle = 1
$borders.Item(-4).LineWidth = 2
$borders.Item(-5).LineStyle = 1
$borders.Item(-5).LineWidth = 2
$borders.Item(-3).LineStyle = 1
$borders.Item(-3).LineWidth = 2
$borders.Item(-6).LineStyle = 1
$borders.Item(-6).LineWidth = 2
